$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value2 = 0.06597369909286499
$ws.Cells.Item(2, 2).Value2 = 0.9862858057022095
$ws.Cells.Item(2, 3).Value2 = 0.1162338629364967
$ws.Cells.Item(2, 4).Value2 = 0.972449004650116
$ws.Cells.Item(3, 1).Value2 = 0.01118937507271767
$ws.Cells.Item(3, 2).Value2 = 0.9980478882789612
$ws.Cells.Item(3, 3).Value2 = 0.1778699904680252
$ws.Cells.Item(3, 4).Value2 = 0.9602678418159485
$ws.Cells.Item(4, 1).Value2 = 0.006237824447453022
$ws.Cells.Item(4, 2).Value2 = 0.9984926581382751
$ws.Cells.Item(4, 3).Value2 = 0.09072687476873398
$ws.Cells.Item(4, 4).Value2 = 0.967793345451355
$ws.Cells.Item(5, 1).Value2 = 0.003342871554195881
$ws.Cells.Item(5, 2).Value2 = 0.9989868998527527
$ws.Cells.Item(5, 3).Value2 = 0.1091148778796196
$ws.Cells.Item(5, 4).Value2 = 0.9642857313156128
$ws.Cells.Item(6, 1).Value2 = 0.001545891398563981
$ws.Cells.Item(6, 2).Value2 = 0.9996787905693054
$ws.Cells.Item(6, 3).Value2 = 0.1252525746822357
$ws.Cells.Item(6, 4).Value2 = 0.9676657915115356
$ws.Cells.Item(7, 1).Value2 = 0.002513718325644732
$ws.Cells.Item(7, 2).Value2 = 0.9992833733558655
$ws.Cells.Item(7, 3).Value2 = 0.1176488399505615
$ws.Cells.Item(7, 4).Value2 = 0.9742984771728516
$ws.Cells.Item(8, 1).Value2 = 0.001212093629874289
$ws.Cells.Item(8, 2).Value2 = 0.9996787905693054
$ws.Cells.Item(8, 3).Value2 = 0.1386571675539017
$ws.Cells.Item(8, 4).Value2 = 0.9737244844436646
$ws.Cells.Item(9, 1).Value2 = 0.001535779680125415
$ws.Cells.Item(9, 2).Value2 = 0.9996293187141418
$ws.Cells.Item(9, 3).Value2 = 0.08428902924060822
$ws.Cells.Item(9, 4).Value2 = 0.9753826260566711
$ws.Cells.Item(10, 1).Value2 = 0.0008829118451103568
$ws.Cells.Item(10, 2).Value2 = 0.99980229139328
$ws.Cells.Item(10, 3).Value2 = 0.1176920011639595
$ws.Cells.Item(10, 4).Value2 = 0.9710459113121033
$ws.Cells.Item(11, 1).Value2 = 0.001806232030503452
$ws.Cells.Item(11, 2).Value2 = 0.9996787905693054
$ws.Cells.Item(11, 3).Value2 = 0.1974737197160721
$ws.Cells.Item(11, 4).Value2 = 0.9653061032295227
$ws.Cells.Item(12, 1).Value2 = 0.0009626789251342416
$ws.Cells.Item(12, 2).Value2 = 0.999777615070343
$ws.Cells.Item(12, 3).Value2 = 0.2835482656955719
$ws.Cells.Item(12, 4).Value2 = 0.9612882733345032
$ws.Cells.Item(13, 1).Value2 = 0.001204178435727954
$ws.Cells.Item(13, 2).Value2 = 0.9998270273208618
$ws.Cells.Item(13, 3).Value2 = 0.2817160487174988
$ws.Cells.Item(13, 4).Value2 = 0.9616709351539612
$ws.Cells.Item(14, 1).Value2 = 0.00106389750726521
$ws.Cells.Item(14, 2).Value2 = 0.9997034668922424
$ws.Cells.Item(14, 3).Value2 = 0.1330789923667908
$ws.Cells.Item(14, 4).Value2 = 0.9748724699020386
$ws.Cells.Item(15, 1).Value2 = 0.0006034984835423529
$ws.Cells.Item(15, 2).Value2 = 0.99980229139328
$ws.Cells.Item(15, 3).Value2 = 0.1247851401567459
$ws.Cells.Item(15, 4).Value2 = 0.9728953838348389
$ws.Cells.Item(16, 1).Value2 = 0.0008754157461225986
$ws.Cells.Item(16, 2).Value2 = 0.99980229139328
$ws.Cells.Item(16, 3).Value2 = 0.1301927864551544
$ws.Cells.Item(16, 4).Value2 = 0.964859664440155
$ws.Cells.Item(17, 1).Value2 = 0.001120057655498385
$ws.Cells.Item(17, 2).Value2 = 0.999777615070343
$ws.Cells.Item(17, 3).Value2 = 0.06684364378452301
$ws.Cells.Item(17, 4).Value2 = 0.9844387769699097
$ws.Cells.Item(18, 1).Value2 = 0.0004083602398168296
$ws.Cells.Item(18, 2).Value2 = 0.9999011754989624
$ws.Cells.Item(18, 3).Value2 = 0.1191778630018234
$ws.Cells.Item(18, 4).Value2 = 0.9788265228271484
$ws.Cells.Item(19, 1).Value2 = 0.0009464031318202615
$ws.Cells.Item(19, 2).Value2 = 0.9998270273208618
$ws.Cells.Item(19, 3).Value2 = 0.04982419684529305
$ws.Cells.Item(19, 4).Value2 = 0.9826530814170837
$ws.Cells.Item(20, 1).Value2 = 0.0003086334909312427
$ws.Cells.Item(20, 2).Value2 = 0.9998517632484436
$ws.Cells.Item(20, 3).Value2 = 0.05738229677081108
$ws.Cells.Item(20, 4).Value2 = 0.983227014541626
$ws.Cells.Item(21, 1).Value2 = 0.0007671194616705179
$ws.Cells.Item(21, 2).Value2 = 0.9998270273208618
$ws.Cells.Item(21, 3).Value2 = 0.05341269075870514
$ws.Cells.Item(21, 4).Value2 = 0.9883928298950195
$ws.Cells.Item(22, 1).Value2 = 0.0002403905964456499
$ws.Cells.Item(22, 2).Value2 = 0.9998764395713806
$ws.Cells.Item(22, 3).Value2 = 0.1418747305870056
$ws.Cells.Item(22, 4).Value2 = 0.9640306234359741
$ws.Cells.Item(23, 1).Value2 = 0.001503307605162263
$ws.Cells.Item(23, 2).Value2 = 0.999777615070343
$ws.Cells.Item(23, 3).Value2 = 0.180696040391922
$ws.Cells.Item(23, 4).Value2 = 0.9615433812141418
$ws.Cells.Item(24, 1).Value2 = 0.000531272089574486
$ws.Cells.Item(24, 2).Value2 = 0.9999011754989624
$ws.Cells.Item(24, 3).Value2 = 0.1591192036867142
$ws.Cells.Item(24, 4).Value2 = 0.9628826379776001
$ws.Cells.Item(25, 1).Value2 = 0.0004924956010654569
$ws.Cells.Item(25, 2).Value2 = 0.999777615070343
$ws.Cells.Item(25, 3).Value2 = 0.1006185412406921
$ws.Cells.Item(25, 4).Value2 = 0.9730229377746582
$ws.Cells.Item(26, 1).Value2 = 0.0001223197468789294
$ws.Cells.Item(26, 2).Value2 = 0.9999752640724182
$ws.Cells.Item(26, 3).Value2 = 0.1592647284269333
$ws.Cells.Item(26, 4).Value2 = 0.9720025658607483
$ws.Cells.Item(27, 1).Value2 = 0.0002229893143521622
$ws.Cells.Item(27, 2).Value2 = 0.9999258518218994
$ws.Cells.Item(27, 3).Value2 = 0.05045899748802185
$ws.Cells.Item(27, 4).Value2 = 0.9850765466690063
$ws.Cells.Item(28, 1).Value2 = 0.001348324003629386
$ws.Cells.Item(28, 2).Value2 = 0.9997528791427612
$ws.Cells.Item(28, 3).Value2 = 0.1724154651165009
$ws.Cells.Item(28, 4).Value2 = 0.9633290767669678
$ws.Cells.Item(29, 1).Value2 = 0.0004411022528074682
$ws.Cells.Item(29, 2).Value2 = 0.9999505877494812
$ws.Cells.Item(29, 3).Value2 = 0.05996633321046829
$ws.Cells.Item(29, 4).Value2 = 0.9857142567634583
$ws.Cells.Item(30, 1).Value2 = 0.0005151263321749866
$ws.Cells.Item(30, 2).Value2 = 0.9999011754989624
$ws.Cells.Item(30, 3).Value2 = 0.08336549997329712
$ws.Cells.Item(30, 4).Value2 = 0.9752551317214966
$ws.Cells.Item(31, 1).Value2 = 0.0004033498698845506
$ws.Cells.Item(31, 2).Value2 = 0.9999011754989624
$ws.Cells.Item(31, 3).Value2 = 0.141542986035347
$ws.Cells.Item(31, 4).Value2 = 0.9673469662666321
$ws.Cells.Item(32, 1).Value2 = 0.0002607712813187391
$ws.Cells.Item(32, 2).Value2 = 0.9998517632484436
$ws.Cells.Item(32, 3).Value2 = 0.2527180016040802
$ws.Cells.Item(32, 4).Value2 = 0.9607780575752258
$ws.Cells.Item(33, 1).Value2 = 0.0006109431269578636
$ws.Cells.Item(33, 2).Value2 = 0.9998517632484436
$ws.Cells.Item(33, 3).Value2 = 0.08328565210103989
$ws.Cells.Item(33, 4).Value2 = 0.9809311032295227
$ws.Cells.Item(34, 1).Value2 = 0.00007301733421627432
$ws.Cells.Item(34, 2).Value2 = 0.9999752640724182
$ws.Cells.Item(34, 3).Value2 = 0.07665637135505676
$ws.Cells.Item(34, 4).Value2 = 0.9818239808082581
$ws.Cells.Item(35, 1).Value2 = 0.0004253085935488343
$ws.Cells.Item(35, 2).Value2 = 0.9998517632484436
$ws.Cells.Item(35, 3).Value2 = 0.06239629536867142
$ws.Cells.Item(35, 4).Value2 = 0.9819515347480774
$ws.Cells.Item(36, 1).Value2 = 0.0003180121420882642
$ws.Cells.Item(36, 2).Value2 = 0.9998764395713806
$ws.Cells.Item(36, 3).Value2 = 0.0878874734044075
$ws.Cells.Item(36, 4).Value2 = 0.9782525300979614
$ws.Cells.Item(37, 1).Value2 = 0.000464294571429491
$ws.Cells.Item(37, 2).Value2 = 0.9999011754989624
$ws.Cells.Item(37, 3).Value2 = 0.2596498727798462
$ws.Cells.Item(37, 4).Value2 = 0.9632652997970581
$ws.Cells.Item(38, 1).Value2 = 0.0005450963508337736
$ws.Cells.Item(38, 2).Value2 = 0.9999011754989624
$ws.Cells.Item(38, 3).Value2 = 0.119734063744545
$ws.Cells.Item(38, 4).Value2 = 0.9732780456542969
$ws.Cells.Item(39, 1).Value2 = 0.00002934718031610828
$ws.Cells.Item(39, 2).Value2 = 1
$ws.Cells.Item(39, 3).Value2 = 0.1474553942680359
$ws.Cells.Item(39, 4).Value2 = 0.969515323638916
$ws.Cells.Item(40, 1).Value2 = 0.0004090869915671647
$ws.Cells.Item(40, 2).Value2 = 0.9998764395713806
$ws.Cells.Item(40, 3).Value2 = 0.06143205985426903
$ws.Cells.Item(40, 4).Value2 = 0.9834821224212646
$ws.Cells.Item(41, 1).Value2 = 0.0001451104908483103
$ws.Cells.Item(41, 2).Value2 = 0.9999505877494812
$ws.Cells.Item(41, 3).Value2 = 0.1485311686992645
$ws.Cells.Item(41, 4).Value2 = 0.9730229377746582
$ws.Cells.Item(42, 1).Value2 = 0.0004318904830142856
$ws.Cells.Item(42, 2).Value2 = 0.9998517632484436
$ws.Cells.Item(42, 3).Value2 = 0.190279483795166
$ws.Cells.Item(42, 4).Value2 = 0.9632652997970581
$ws.Cells.Item(43, 1).Value2 = 0.0001521644444437698
$ws.Cells.Item(43, 2).Value2 = 0.9999752640724182
$ws.Cells.Item(43, 3).Value2 = 0.1031414419412613
$ws.Cells.Item(43, 4).Value2 = 0.9774234890937805
$ws.Cells.Item(44, 1).Value2 = 0.00009628408588469028
$ws.Cells.Item(44, 2).Value2 = 0.9999752640724182
$ws.Cells.Item(44, 3).Value2 = 0.4489351809024811
$ws.Cells.Item(44, 4).Value2 = 0.9612244963645935
$ws.Cells.Item(45, 1).Value2 = 0.00009661566582508385
$ws.Cells.Item(45, 2).Value2 = 0.9999752640724182
$ws.Cells.Item(45, 3).Value2 = 0.1094375625252724
$ws.Cells.Item(45, 4).Value2 = 0.9802296161651611
$ws.Cells.Item(46, 1).Value2 = 0.0004403672646731138
$ws.Cells.Item(46, 2).Value2 = 0.9999505877494812
$ws.Cells.Item(46, 3).Value2 = 0.1585179269313812
$ws.Cells.Item(46, 4).Value2 = 0.9688137769699097
$ws.Cells.Item(47, 1).Value2 = 0.0006111921975389123
$ws.Cells.Item(47, 2).Value2 = 0.9999011754989624
$ws.Cells.Item(47, 3).Value2 = 0.131764829158783
$ws.Cells.Item(47, 4).Value2 = 0.9688775539398193
$ws.Cells.Item(48, 1).Value2 = 0.0006102813640609384
$ws.Cells.Item(48, 2).Value2 = 0.9999258518218994
$ws.Cells.Item(48, 3).Value2 = 0.2706469595432281
$ws.Cells.Item(48, 4).Value2 = 0.965050995349884
$ws.Cells.Item(49, 1).Value2 = 0.0003370883350726217
$ws.Cells.Item(49, 2).Value2 = 0.9998517632484436
$ws.Cells.Item(49, 3).Value2 = 0.1378201246261597
$ws.Cells.Item(49, 4).Value2 = 0.9745535850524902
$ws.Cells.Item(50, 1).Value2 = 0.0001446453388780355
$ws.Cells.Item(50, 2).Value2 = 0.9999505877494812
$ws.Cells.Item(50, 3).Value2 = 0.1958103626966476
$ws.Cells.Item(50, 4).Value2 = 0.9702805876731873
$ws.Cells.Item(51, 1).Value2 = 0.0002668813685886562
$ws.Cells.Item(51, 2).Value2 = 0.9999011754989624
$ws.Cells.Item(51, 3).Value2 = 0.1921824812889099
$ws.Cells.Item(51, 4).Value2 = 0.9708545804023743
